$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Oppsett for programvariabel (A2): bytt årstall 2021 -> 2022 ---
$ws.Range("A2").Value = "Studiebarometeret 2022"

# --- Ny rad 60: progresjonsvariabel (Andel heltid / progresjon) ---
$ws.Range("C60").Value = "Andel heltid"
$ws.Range("C60").Font.Bold = $false
$ws.Range("C60").HorizontalAlignment = -4131
$ws.Range("C60").VerticalAlignment = -4160

$ws.Range("D60").Value = "progresjon"
$ws.Range("D60").Font.Bold = $false
$ws.Range("D60").HorizontalAlignment = -4131
$ws.Range("D60").VerticalAlignment = -4160

$ws.Range("E60").Value = $false
$ws.Range("E60").Font.Bold = $false
$ws.Range("E60").HorizontalAlignment = -4131
$ws.Range("E60").VerticalAlignment = -4160

# --- Ny kolonne F: Format, for å fargekode tidsvariabelen ---
$ws.Range("F1").Value = "Format"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4131
$ws.Range("F1").VerticalAlignment = -4160

$ws.Range("F60").Value = "prosent"
$ws.Range("F60").Font.Bold = $false
$ws.Range("F60").HorizontalAlignment = -4131
$ws.Range("F60").VerticalAlignment = -4160

# --- Oppdater visning: fjern gammel rulle-/markeringsposisjon, still til F1 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F1").Select()
